$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "abu10@gmail.com"
$ws.Range("B3").Value = "pallu10@gmail.com"
$ws.Range("B4").Value = "zarina10@gmail.com"
$ws.Range("B5").Value = "tahira10@gmail.com"

$ws.Range("C5").Select()
